# agrego número de orden de pago en la validación de pagos en SISE
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5: replace the old siniestro number with a new one, and bump the importe
$ws.Range("B5").Value = "'0420172008483 "
$ws.Range("C5").Value = 120

# New row 6
$ws.Range("A6").Value = "PREPROD"
$ws.Range("B6").Value = "'1220170301429 "
$ws.Range("C6").Value = 100

# New header in D1
$ws.Range("D1").Value = "NumOrden"

# New row 7
$ws.Range("A7").Value = "PREPROD"
$ws.Range("B7").Value = "'1120170200937 "
$ws.Range("C7").Value = 100
$ws.Range("D7").Value = "'1703271"

# Selection moves to B8 per the saved view state
$ws.Range("B8").Select()
